$wb = $excel.ActiveWorkbook

# The "model" sheet (3rd tab) gains a new first column with extracted
# eigenvalue / EFA output-file info used for the EFA screeplot.
$ws = $wb.Worksheets.Item(3)

$ws.Columns.Item(1).Insert()

$ws.Range("A1").Value = "model_efa"
$ws.Range("A2").Value = "CS1_EFA.out"
$ws.Range("A3").Value = "CS2_EFA.out"
$ws.Range("A4").Value = "CS3_EFA.out"

# Make "model" the active sheet/tab, with A5 selected, instead of "text".
$ws.Activate()
$ws.Range("A5").Select()
